$wb = $excel.ActiveWorkbook

# --- Sprint sheet: mark the remaining tasks (rows 7-14) as "completed" ---
$sprint = $wb.Worksheets.Item("Sprint")
$sprint.Range("C7:C14").Value = "completed"

# --- Move the active tab / selection from "Burndown Chart" to "Sprint" ---
$sprint.Activate()
$sprint.Range("C15").Select()
